$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.765.82"
$ws.Range("E2").Value = "  -2.37%  "
$ws.Range("D3").Value = "3.086.01"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'525.64"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "'140.86"
$ws.Range("E6").Value = "  -2.58%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.084.54"
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  -3.81%  "
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").Value = "3.615.65"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").Value = "'25.40"
$ws.Range("E15").Value = "  -6.49%  "
$ws.Range("E16").Value = "  -2.06%  "
$ws.Range("D17").Value = "57.810.06"
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").Value = "3.086.72"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "'6.07"
$ws.Range("E19").Value = "  -2.70%  "
$ws.Range("E20").Value = "  -3.54%  "
$ws.Range("E21").Value = "  -3.84%  "
$ws.Range("D22").Value = "'341.18"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'67.22"
$ws.Range("E25").Value = "  +1.70%  "
$ws.Range("E26").Value = "  -2.10%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "0.0₃0911"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D30").Value = "'6.38"
$ws.Range("E30").Value = "  -4.37%  "
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("D33").Value = "'20.91"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("E34").Value = "  -4.20%  "
$ws.Range("D35").Value = "'158.74"
$ws.Range("E35").Value = "  +2.02%  "
$ws.Range("D36").Value = "'4.60"
$ws.Range("E36").Value = "  -1.54%  "
$ws.Range("D37").Value = "'6.13"
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("D38").Value = "'25.87"
$ws.Range("E38").Value = "  -5.45%  "
$ws.Range("E39").Value = "  -4.86%  "
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("D41").Value = "'1.57"
$ws.Range("E41").Value = "  +7.39%  "
$ws.Range("D42").Value = "'3.98"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").Value = "'0.682"
$ws.Range("E43").Value = "  +2.49%  "
$ws.Range("D44").Value = "3.125.42"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("D45").Value = "'36.88"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "'0.0261"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("D48").Value = "2.273.50"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("E49").Value = "  +2.73%  "
$ws.Range("D50").Value = "'6.08"
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").Value = "'20.46"
$ws.Range("E51").Value = "  -3.24%  "
